$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: membership count. Force as literal text "1,500" (matches original
# cell which stored "253" as a shared string / text, not a number), then
# strip the style Excel auto-applies for a numeric-looking text entry so
# the cell keeps the workbook's default (unstyled) formatting.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1,500"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "No, Soleo Health does not encompass community sites.`nJustification: Soleo Health is a specialty pharmacy that focuses on providing infusion therapy services and support to patients outside of the traditional healthcare setting."
$ws.Range("D2").Value = "No, Soleo Health is not influential on state or local policy. This society primarily focuses on healthcare services rather than policy advocacy."
$ws.Range("E2").Value = "Yes, Soleo Health provides engagement opportunities with leadership. The society offers various channels for members to interact and engage with its leadership team, fostering a sense of community and transparency."
$ws.Range("F2").Value = "No, Soleo Health does not provide support for clinical trial recruitment. Soleo Health is a specialty pharmacy that focuses on providing infusion therapies and nursing services, not on conducting clinical trials."
$ws.Range("G2").Value = "No, Soleo Health does not provide engagement opportunities with payors. Soleo Health is focused on providing specialty pharmacy and infusion services rather than directly engaging with payors."
$ws.Range("H2").Value = "No, there is no information available about area experts on the board of Soleo Health."
$ws.Range("I2").Value = "No, Soleo Health, as a specialty pharmacy, focuses on providing care and services related to the delivery of medications for complex and chronic conditions, rather than being directly involved in therapeutic research collaborations."
$ws.Range("J2").Value = "No, Soleo Health does not include top therapeutic area experts on its board. Soleo Health focuses on providing innovative specialty pharmacy and infusion services."
$ws.Range("L2").Value = "2025-03-13 09:14:18"

# C2 contains an embedded line break, which makes Excel auto-expand the row
# height. Re-run AutoFit so the row's height reverts to the sheet default
# (no explicit ht/customHeight override), matching the rest of the workbook.
$ws.Rows(2).AutoFit()
